$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.27%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "29.50"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.61%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.111"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "3.45%"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.18%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.372"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "1.78%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.409"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.48%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.365"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.41%"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "1.09%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1590"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2.24%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06713"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-1.56%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07703"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-0.27%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.02950"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-1.24%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.22%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001572"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "-1.47%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04512"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.84%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0006449"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-1.31%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006277"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.95%"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-0.44%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.220"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-1.08%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3216"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "2.09%"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-2.93%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.082"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.31%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.84%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001190"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "0.10%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-4.55%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "5.76%"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "-1.22%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04226"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "1.45%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006730"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "1.49%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1241"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "0.55%"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "-7.62%"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-2.70%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005701"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "2.57%"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.970"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "26.11%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-29.46%"
